$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5252051949501038
$ws.Range("B1").Value = 1.647102952003479
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.40656566619873
$ws.Range("E1").Value = 1.386918067932129
